# Edit script for C1--C2-and-C3-PowerPoint.pptx
#
# Two logical changes (per the authoritative OOXML diff):
#
#  1. On slide 16, the table's table-style reference changes from
#     {9A31AE2E-AA73-486C-AD61-B1066D3E96F1} to
#     {9463B974-5159-407A-BD90-BD414A94BBB9}.
#
#  2. The deck's theme colour values are swapped: the theme that is
#     actually bound to the slide master / presentation (the "Integral"
#     palette) is recoloured to the stock "Office Theme" palette.
#     (dk1/lt1 are pure black/white in both palettes, so only the other
#     ten slots actually change value.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------

$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{9463B974-5159-407A-BD90-BD414A94BBB9}")
    }
}

# --- 2. Theme colour swap ---------------------------------------------------

$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
